$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list - Price (D) and Volume(1h) (E) columns
# Cells whose new D value is a plain decimal number need NumberFormat "@"
# forced first, otherwise Excel auto-converts the text to a numeric value
# (e.g. "139.30" -> 139.3) and drops the trailing zero / string type.

$ws.Range("D2").Value = '60.644.87'
$ws.Range("E2").Value = '  -2.24%  '
$ws.Range("D3").Value = '2.398.12'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.25'
$ws.Range("E5").Value = '  -1.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.30'
$ws.Range("E6").Value = '  -2.50%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.525'
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("D9").Value = '2.378.38'
$ws.Range("E9").Value = '  -2.46%  '
$ws.Range("E10").Value = '  +0.94%  '
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("E12").Value = '  -2.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.339'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.82'
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000169'
$ws.Range("E15").Value = '  -1.23%  '
$ws.Range("D16").Value = '2.826.25'
$ws.Range("E16").Value = '  -0.94%  '
$ws.Range("D17").Value = '60.654.34'
$ws.Range("E17").Value = '  -2.19%  '
$ws.Range("D18").Value = '2.374.94'
$ws.Range("E18").Value = '  -2.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.53'
$ws.Range("E19").Value = '  -2.53%  '
$ws.Range("E20").Value = '  +1.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.84'
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.10'
$ws.Range("E23").Value = '  +2.02%  '
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.87'
$ws.Range("E25").Value = '  -5.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.56'
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("E27").Value = '  -8.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '570.65'
$ws.Range("E28").Value = '  -6.20%  '
$ws.Range("D29").Value = '2.504.24'
$ws.Range("E29").Value = '  -2.40%  '
$ws.Range("D30").Value = '0.0₃0910'
$ws.Range("E30").Value = '  -3.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.84'
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("E32").Value = '  -4.99%  '
$ws.Range("E33").Value = '  -2.47%  '
$ws.Range("E34").Value = '  -6.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.60'
$ws.Range("E36").Value = '  -5.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.367'
$ws.Range("E37").Value = '  -2.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.37'
$ws.Range("E38").Value = '  -3.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '147.60'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.13'
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("E41").Value = '  -3.98%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  -3.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.78'
$ws.Range("E44").Value = '  -4.09%  '
$ws.Range("E45").Value = '  -4.12%  '
$ws.Range("D46").Value = '0.0₆0279'
$ws.Range("E46").Value = '  +19.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '139.30'
$ws.Range("E47").Value = '  -2.21%  '
$ws.Range("E48").Value = '  -3.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.583'
$ws.Range("E49").Value = '  -3.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0502'
$ws.Range("E50").Value = '  -3.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.25'
$ws.Range("E51").Value = '  -0.68%  '
